$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# New Product Backlog items added for the "room/level/building count",
# "window area" and "rent" functionality. The acceptance-criteria column (E)
# was filled in first, followed by the backlog-item column (A), which is the
# order the new shared strings were originally recorded in.
$ws.Range("E11").Value = "o Dla pomieszczenie zwraca 1`no Dla poziomu wylicza ilość pomieszczeń`no Dla budynku wyliczane jako suma na poziomach"
$ws.Range("E12").Value = "o Prawidłowo podawane dla pomieszczenia`no Dla poziomu wyliczane sumę powierzchni okien`no Dla budynku wyliczane jest jako suma powierzchni okien na poziomach"
$ws.Range("E13").Value = "o Prawidłowo podawane dla pomieszczenia`no Dla poziomu wyliczane sumę czynszów`no Dla budynku wyliczane jest jako suma czynszów poziomów"
$ws.Range("E14").Value = "o Zwroci liste pokoi o nie wyższym czynszu niż zadany"
$ws.Range("E15").Value = "o Zwroci liste pokoi o nie mniejszej powierzchni niż zadana"

$ws.Range("A11").Value = "Jako zarządca budynku mogę sprawdzić liczbę pomieszczeń na danym piętrze lub w całym budynku"
$ws.Range("A12").Value = "Jako zarządca budynku mogę sprawdzić powierzchnię okien w pomieszczeniu, na danym piętrze i w całym buynku"
$ws.Range("A13").Value = "Jako zarządca budynku mogę sprawdzić czynsz za wynajem pokoju, poziomu lub budynku"
$ws.Range("A14").Value = "Jako zarządca budynku mogę sprzawdzić listę pokoi o czynszach mniejszych lub rownych niż zadana wartosc"
$ws.Range("A15").Value = "Jako zarządca budynku mogę sprawdzić listę pokoi o powierzchni większej lub równej od zadanej wartości"

# The edited workbook was left with the "Product Backlog" sheet active/selected.
[void]$ws.Range("A15").Select()
